$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "2 0 0"

# Row 2 updates (reuses shared string index 4 for B2)
$ws.Cells.Item(2,2).Value = "2023-02-02T00:00"
$ws.Cells.Item(2,3).Value = 3.0
$ws.Cells.Item(2,4).Value = 1500.0

# Row 3 (new) - introduces shared string index 5
$ws.Cells.Item(3,1).Value = 1.0
$ws.Cells.Item(3,2).Value = "2024-03-02T00:00"
$ws.Cells.Item(3,3).Value = 3.0
$ws.Cells.Item(3,4).Value = 1500.0

# Row 4 (new) - introduces shared string index 6
$ws.Cells.Item(4,1).Value = 1.0
$ws.Cells.Item(4,2).Value = "2024-02-02T00:00"
$ws.Cells.Item(4,3).Value = 2.0
$ws.Cells.Item(4,4).Value = 1000.0
